# Update countries & provincias Spain
# - Refresh COVID numbers for a handful of countries
# - Re-sort the country table by "Casos totales" (column B) descending,
#   since the updated totals change the ranking
# - Bump the "last updated" timestamp in A1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New per-country figures: Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes
$updates = @{
    "Estados Unidos" = @(3770138, 126, 1741398, 1886675, 0, 1, 142065)
    "India"          = @(1040948, 491, 654193, 360460, 0, 10, 26295)
    "Rusia"          = @(765437, 6234, 546863, 206327, 0, 124, 12247)
    "Filipinas"      = @(65304, 2357, 22067, 41464, 0, 113, 1773)
    "Ucrania"        = @(58111, 847, 30525, 26109, 0, 21, 1477)
    "Singapur"       = @(47655, 202, 43577, 4051, 0, 0, 27)
    "Hungria"        = @(4315, 22, 3222, 497, 0, 1, 596)
    "Estonia"        = @(2021, 1, 1912, 40, 0, 0, 69)
    "Montenegro"     = @(2024, 59, 378, 1618, 0, 2, 28)
    "Eslovaquia"     = @(1976, 11, 1523, 425, 0, 0, 28)
    "Letonia"        = @(1189, 4, 1022, 136, 0, 0, 31)
    "Gambia"         = @(93, 15, 49, 40, 0, 1, 4)
}

$lookupRange = $ws.Range("A4:A219")

foreach ($country in $updates.Keys) {
    $values = $updates[$country]
    $hit = $lookupRange.Find($country)
    if ($hit -eq $null) {
        continue
    }
    $row = $hit.Row()
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
    $ws.Cells.Item($row, 8).Value = $values[6]
}

# Re-rank the table now that totals changed
$dataRange = $ws.Range("A4:H219")
$keyRange = $ws.Range("B4:B219")
$dataRange.Sort($keyRange, 2)

# Bump the timestamp shown above the table
$ws.Range("A1").Value = "Datos actualizados a 18 de Julio de 2020 a las 10:10"
